$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A: code names for rows 35-39
$ws.Range("A35").Value = "tavg"
$ws.Range("A36").Value = "tmin"
$ws.Range("A37").Value = "tmax"
$ws.Range("A38").Value = "prcp"
$ws.Range("A39").Value = "wspd"

# Column B: column numbers for rows 35-39
$ws.Range("B35").Value = 34
$ws.Range("B36").Value = 35
$ws.Range("B37").Value = 36
$ws.Range("B38").Value = 37
$ws.Range("B39").Value = 38

# Column C: data type for rows 35-39
$ws.Range("C35").Value = "Float64"
$ws.Range("C36").Value = "Float64"
$ws.Range("C37").Value = "Float64"
$ws.Range("C38").Value = "Float64"
$ws.Range("C39").Value = "Float64"

# Row 40 (pres) added afterwards
$ws.Range("A40").Value = "pres"
$ws.Range("B40").Value = 39
$ws.Range("C40").Value = "Float64"

# Column D: units for rows 35-40
$ws.Range("D35").Value = "°C"
$ws.Range("D36").Value = "°C"
$ws.Range("D37").Value = "°C"
$ws.Range("D38").Value = "mm"
$ws.Range("D39").Value = "km/h"
$ws.Range("D40").Value = "hPa"

# Column E: group for rows 35-40
$ws.Range("E35").Value = "Weather"
$ws.Range("E36").Value = "Weather"
$ws.Range("E37").Value = "Weather"
$ws.Range("E38").Value = "Weather"
$ws.Range("E39").Value = "Weather"
$ws.Range("E40").Value = "Weather"

# Column F: descriptions for rows 35-40
$ws.Range("F35").Value = "Average air temperature"
$ws.Range("F36").Value = "Minimum air temperature"
$ws.Range("F37").Value = "Maximum air temperature"
$ws.Range("F38").Value = "Precipitation total"
$ws.Range("F39").Value = "Wind speed"
$ws.Range("F40").Value = "Average sea-level air pressure"

$ws.Range("F40").Select()
